# Increment the "Förändrad" (Changed) date in column C by one day
# for every data row (rows 2 through 536), i.e. 45214 -> 45215.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) {
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = $current + 1
    }
}
